$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Semester"
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 3

$ws.Range("E7").Select()
